$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing ticket summary ("Login") for the new ticket row (row 3, column B)
$ws.Range("B3").Value = "Login"

# Move the active selection to reflect where the user left off editing
$ws.Range("B5").Select()
